# AfDD_2023_Annex_Table_Tab31.xlsx -- "Add files via upload" re-export.
#
# The authoritative diff for this commit consists almost entirely of
# re-save noise produced when the workbook was opened and re-saved by a
# newer Excel build (rupBuild 25601 -> 26130): fresh random revision and
# window GUIDs, a bumped fileVersion/calcPr, window geometry, and tiny
# font-metric-driven shifts in default row height / bestFit column
# widths / x14ac:dyDescent that Excel recomputes on save. None of that
# is exposed as document content through the Excel object model, so it
# cannot (and should not) be faked here.
#
# The one substantive, content-level change is a silent recalculation
# of the cached F ("gross" index) / G ("net" index) values on every
# continent-aggregate row of the table: the "Africa" rows (13/23/38/
# 45/61) plus the repeated regional summary block (rows 62-98). Those
# numbers differ only in the 8th-10th significant digit (recalculated
# upstream -- no formulas are stored in this sheet, every cell is a
# pasted/cached constant), so we simply replay the new cached values
# cell by cell exactly as they appear in the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab31")
$ws.Range("F13").Value = 28.9458738496504
$ws.Range("G13").Value = 33.239513643544903
$ws.Range("F23").Value = 35.551194138485997
$ws.Range("G23").Value = 39.141113217542497
$ws.Range("F38").Value = 27.8583101456707
$ws.Range("G38").Value = 37.024542271840403
$ws.Range("F45").Value = 29.187883490849298
$ws.Range("G45").Value = 32.2716211970004
$ws.Range("F61").Value = 45.045905757450598
$ws.Range("G61").Value = 50.063165609684901
$ws.Range("F62").Value = 33.097542292498098
$ws.Range("G62").Value = 38.084196775311398
$ws.Range("F63").Value = 16.667448396873901
$ws.Range("G63").Value = 18.8561065918148
$ws.Range("F64").Value = 29.830368191411399
$ws.Range("G64").Value = 35.156048141931699
$ws.Range("F65").Value = 14.622214384119699
$ws.Range("G65").Value = 17.675093100166901
$ws.Range("F66").Value = 17.4503095022534
$ws.Range("G66").Value = 19.796419162071999
$ws.Range("F67").Value = 29.506919722778001
$ws.Range("G67").Value = 33.914813030885803
$ws.Range("F68").Value = 35.926316842943798
$ws.Range("G68").Value = 41.186639592389298
$ws.Range("F69").Value = 34.396623528155096
$ws.Range("G69").Value = 39.649967043318597
$ws.Range("F70").Value = 37.485597620657501
$ws.Range("G70").Value = 40.566468118392997
$ws.Range("F71").Value = 45.045905757450598
$ws.Range("G71").Value = 50.063165609684901
$ws.Range("F72").Value = 24.4613588914941
$ws.Range("G72").Value = 33.170413475191303
$ws.Range("F73").Value = 31.390991222289099
$ws.Range("G73").Value = 35.746579179312199
$ws.Range("F74").Value = 29.076790381740999
$ws.Range("G74").Value = 33.1571739322086
$ws.Range("F75").Value = 37.194096678801898
$ws.Range("G75").Value = 42.327929194931599
$ws.Range("F76").Value = 23.185482896023601
$ws.Range("G76").Value = 26.541721292280599
$ws.Range("F77").Value = 30.7155199519975
$ws.Range("G77").Value = 36.5236278875389
$ws.Range("F78").Value = 18.529082723174401
$ws.Range("G78").Value = 19.552133054814899
$ws.Range("F79").Value = 14.413877959767699
$ws.Range("G79").Value = 15.3533028677018
$ws.Range("F80").Value = 31.373435624011002
$ws.Range("G80").Value = 36.021318858723198
$ws.Range("F81").Value = 28.327735479727199
$ws.Range("G81").Value = 31.128155317501399
$ws.Range("F82").Value = 33.332031878584303
$ws.Range("G82").Value = 38.413241384598898
$ws.Range("F83").Value = 15.6443961640512
$ws.Range("G83").Value = 17.778334159390301
$ws.Range("F84").Value = 29.303801327769801
$ws.Range("G84").Value = 38.722121784835402
$ws.Range("F86").Value = 35.747517889135899
$ws.Range("G86").Value = 39.984263488514102
$ws.Range("F87").Value = 20.067648962232699
$ws.Range("G87").Value = 23.742132608217801
$ws.Range("G88").Value = 29.030152891861398
$ws.Range("F89").Value = 19.7505540214718
$ws.Range("G89").Value = 23.046088110990699
$ws.Range("F90").Value = 13.4701599887384
$ws.Range("G90").Value = 14.398615713672299
$ws.Range("F91").Value = 32.758286571229497
$ws.Range("G91").Value = 41.347613258315299
$ws.Range("F92").Value = 26.992784866973199
$ws.Range("G92").Value = 37.453315016141197
$ws.Range("G93").Value = 25.183442234125199
$ws.Range("F94").Value = 16.300299272693199
$ws.Range("G94").Value = 20.214241052790999
$ws.Range("F95").Value = 32.216467360322603
$ws.Range("G95").Value = 38.544256091100699
$ws.Range("F96").Value = 38.473727070833299
$ws.Range("G96").Value = 39.630360849614199
$ws.Range("F97").Value = 39.780727852131903
$ws.Range("G97").Value = 45.8455662428415
$ws.Range("F98").Value = 26.235035491502199
$ws.Range("G98").Value = 31.150406540272101
